$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix "Objetivos:" row (row 10): replace placeholder text with correct Portuguese objectives paragraph
$ws.Range("B10").Value = 'a) Fornecer ao aluno uma visão integrada do processo de desenvolvimento de produto, desde as etapas iniciais de geração da ideia, avaliação econômica e desenvolvimento do conceito do produto até a preparação da fábrica, produção e lançamento do produto. b) Apresentar como os principais conhecimentos da formação de engenheiro de produção podem ser aplicados no processo de desenvolvimento de produto. c) Propiciar uma experiência prática de projeto.'
$ws.Range("C10").Value = 'a) Fornecer ao aluno uma visão integrada do processo de desenvolvimento de produto, desde as etapas iniciais de geração da ideia, avaliação econômica e desenvolvimento do conceito do produto até a preparação da fábrica, produção e lançamento do produto. b) Apresentar como os principais conhecimentos da formação de engenheiro de produção podem ser aplicados no processo de desenvolvimento de produto. c) Propiciar uma experiência prática de projeto.'

# 2. Insert a new row at 13 to hold the "Docentes responsaveis" value (professor name),
#    which shifts old rows 13-23 down to 14-24.
$ws.Rows.Item(13).Insert()

# The newly inserted row only carries column A formatting (copied from row 12 above, which has only column A).
# Remove the stray A13 cell entirely so the row matches the target (only B/C populated).
$ws.Range("A13").Clear()

# Copy the B/C formatting (styles) from the row now holding "Programa resumido:" (row 14) onto B13:C13
$ws.Range("B14:C14").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("B13").Value = "5701460 - Antonio Iacono"
$ws.Range("C13").Value = "5701460 - Antonio Iacono"

# 3. Fix "Programa resumido:" row (now row 14): replace "Semestral" with the full Portuguese summary paragraph
$ws.Range("B14").Value = '1. Conceito de processo de negócio e modelo de referência; 2. Desenvolvimento de produto como um processo: Fases principais do processo de desenvolvimento de produto; 3. Definição e conceitos básicos de gerenciamento de projetos; 4. Planejamento estratégico do produto; 5. Planejamento do Projeto; 6. Projeto Informacional; 7. Projeto Conceitual; 8. Projeto Detalhado; 9. Lançamento e Pós-Desenvolvimento. 10. Empreendedorismo e processo de desenvolvimento de produto.'
$ws.Range("C14").Value = '1. Conceito de processo de negócio e modelo de referência; 2. Desenvolvimento de produto como um processo: Fases principais do processo de desenvolvimento de produto; 3. Definição e conceitos básicos de gerenciamento de projetos; 4. Planejamento estratégico do produto; 5. Planejamento do Projeto; 6. Projeto Informacional; 7. Projeto Conceitual; 8. Projeto Detalhado; 9. Lançamento e Pós-Desenvolvimento. 10. Empreendedorismo e processo de desenvolvimento de produto.'

# 4. Fix "Programa:" row (now row 16): replace the leaked date with the full Portuguese program paragraph
$ws.Range("B16").Value = '1. Visão geral do processo de desenvolvimento de produto. Conceito de processo de negócio. Modelo de referência e desenvolvimento de produto como um processo. Fases principais do processo de desenvolvimento de produto. Definição e conceitos básicos de gerenciamento de projetos. 2. Apresentação de um modelo de referência para desenvolvimento de produto. Descrição do modelo de referência. Importância e uso do modelo de referência. Principais características do modelo de referência para desenvolvimento de produto. 3. Apresentação das fases do processo de desenvolvimento de produto. Nesta parte são apresentados os conceitos e metodologias principais de cada fase do desenvolvimento de produto. 3.1. Planejamento estratégico de produtos. Conteúdo sobre a importância e como preparar a estratégia de produto, realizar o gerenciamento da carteira de projetos e elaborar a Minuta do Projeto. 3.2. Planejamento do projeto. Conteúdo sobre técnicas de planejamento do projeto e análise de viabilidade econômica. 3.3. Projeto informacional. Conceitos gerais de pesquisa de mercado, tipos de fontes de dados para pesquisa de mercado, diferenças entre as fontes de dados, formas de utilização de pesquisas e aprofundamento no tipo de pesquisa baseado em enquete. 3.4. Projeto conceitual. Apresentam-se os diversos métodos criativos para concepção, tais como Diagrama de funções, Matriz Morfológica e Matriz de Pugh, e demonstra-se com exemplos como utilizá-los de forma a gerar e escolher uma concepção de produto robusta.3.5. Projeto detalhado e preparação da produção. Apresenta-se o objetivo e resumo das atividades realizadas nas duas fases. No pós-desenvolvimento discute-se especialmente os aspectos relacionados com o meio ambiente, como descarte e reuso de peças do produto. 4. Empreendedorismo e processo de desenvolvimento de produto. O ambiente de inovação e o processo de desenvolvimento de produtos: empresas nascentes de base tecnológica, startups e aceleradoras.'
$ws.Range("C16").Value = '1. Visão geral do processo de desenvolvimento de produto. Conceito de processo de negócio. Modelo de referência e desenvolvimento de produto como um processo. Fases principais do processo de desenvolvimento de produto. Definição e conceitos básicos de gerenciamento de projetos. 2. Apresentação de um modelo de referência para desenvolvimento de produto. Descrição do modelo de referência. Importância e uso do modelo de referência. Principais características do modelo de referência para desenvolvimento de produto. 3. Apresentação das fases do processo de desenvolvimento de produto. Nesta parte são apresentados os conceitos e metodologias principais de cada fase do desenvolvimento de produto. 3.1. Planejamento estratégico de produtos. Conteúdo sobre a importância e como preparar a estratégia de produto, realizar o gerenciamento da carteira de projetos e elaborar a Minuta do Projeto. 3.2. Planejamento do projeto. Conteúdo sobre técnicas de planejamento do projeto e análise de viabilidade econômica. 3.3. Projeto informacional. Conceitos gerais de pesquisa de mercado, tipos de fontes de dados para pesquisa de mercado, diferenças entre as fontes de dados, formas de utilização de pesquisas e aprofundamento no tipo de pesquisa baseado em enquete. 3.4. Projeto conceitual. Apresentam-se os diversos métodos criativos para concepção, tais como Diagrama de funções, Matriz Morfológica e Matriz de Pugh, e demonstra-se com exemplos como utilizá-los de forma a gerar e escolher uma concepção de produto robusta.3.5. Projeto detalhado e preparação da produção. Apresenta-se o objetivo e resumo das atividades realizadas nas duas fases. No pós-desenvolvimento discute-se especialmente os aspectos relacionados com o meio ambiente, como descarte e reuso de peças do produto. 4. Empreendedorismo e processo de desenvolvimento de produto. O ambiente de inovação e o processo de desenvolvimento de produtos: empresas nascentes de base tecnológica, startups e aceleradoras.'

# 5. Fix "Metodo:" row (now row 19): replace the misplaced professor name with "Provas e trabalho"
$ws.Range("B19").Value = "Provas e trabalho"
$ws.Range("C19").Value = "Provas e trabalho"

# 6. Fix "Criterio:" row (now row 20): replace "Provas e trabalho" with the grading formula text
$ws.Range("B20").Value = 'M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre desenvolvimento de produtosM = Média de aproveitamento do aluno'
$ws.Range("C20").Value = 'M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre desenvolvimento de produtosM = Média de aproveitamento do aluno'

# 7. Fix "Norma de recuperacao:" row (now row 21): replace the grading formula text with the recovery-norm text
$ws.Range("B21").Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$ws.Range("C21").Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'

# 8. Fix "Bibliografia:" row (now row 22): replace the recovery-norm text with the bibliography paragraph
$ws.Range("B22").Value = 'ROZENFELD, H.; FORCELLINI, F.A.; AMARAL, D.C.; TOLEDO, J. A.; SILVA, S.L.; ALLIPRANDINI, D.H.; SCALICE, R.K. Gestão de desenvolvimento de produtos: uma referência para a melhoria do processo. São Paulo: Saraiva, 2006. 542p. CHENG, L.C. QFD: desdobramento da função qualidade na gestão de desenvolvimento de produtos. São Paulo: Edgard Blücher, 2007. DENIS, R. C. Uma introdução à história do design. São Paulo: Edgard Blücher., 2000. MALHOTRA, N.K. Pesquisa de marketing: uma orientação aplicada. Porto Alegre: Bookman, 2006. PAHL,G.; BEITZ,W.; FELDHUSEN,J.;GROTE,K. Projeto na engenharia. São Paulo: Edgard Blucher, 2005.'
$ws.Range("C22").Value = 'ROZENFELD, H.; FORCELLINI, F.A.; AMARAL, D.C.; TOLEDO, J. A.; SILVA, S.L.; ALLIPRANDINI, D.H.; SCALICE, R.K. Gestão de desenvolvimento de produtos: uma referência para a melhoria do processo. São Paulo: Saraiva, 2006. 542p. CHENG, L.C. QFD: desdobramento da função qualidade na gestão de desenvolvimento de produtos. São Paulo: Edgard Blücher, 2007. DENIS, R. C. Uma introdução à história do design. São Paulo: Edgard Blücher., 2000. MALHOTRA, N.K. Pesquisa de marketing: uma orientação aplicada. Porto Alegre: Bookman, 2006. PAHL,G.; BEITZ,W.; FELDHUSEN,J.;GROTE,K. Projeto na engenharia. São Paulo: Edgard Blucher, 2005.'
